$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'23.262.38"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.04%  '

$ws.Range('D3').Value = "'1.610.25"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.10%  '

$ws.Range('D4').Value = "'0.9992"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.33%  '

$ws.Range('D5').Value = "'0.9998"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.16%  '

$ws.Range('D6').Value = "'305.70"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.30%  '

$ws.Range('D7').Value = "'0.3760"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.28%  '

$ws.Range('D8').Value = "'53.08"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.46%  '

$ws.Range('E9').Value = '  +0.07%  '

$ws.Range('D10').Value = "'1.265"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.69%  '

$ws.Range('D11').Value = "'0.08150"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.60%  '

$ws.Range('D12').Value = "'0.9994"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.32%  '

$ws.Range('D13').Value = "'22.91"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.46%  '

$ws.Range('D14').Value = "'6.612"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.53%  '

$ws.Range('D15').Value = "'7.372"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.30%  '

$ws.Range('D16').Value = "'0.00001247"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.22%  '

$ws.Range('D17').Value = "'1.609.07"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.48%  '

$ws.Range('D18').Value = "'94.30"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.17%  '

$ws.Range('D19').Value = "'0.06919"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.31%  '

$ws.Range('D20').Value = "'18.23"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.11%  '

$ws.Range('D21').Value = "'6.551"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.18%  '

$ws.Range('E22').Value = '  -0.02%  '

$ws.Range('D23').Value = "'12.92"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.04%  '

$ws.Range('D24').Value = "'23.257.77"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.05%  '

$ws.Range('D25').Value = "'3.098"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +10.49%  '

$ws.Range('D26').Value = "'2.407"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.59%  '

$ws.Range('E27').Value = '  +1.22%  '

$ws.Range('D28').Value = "'150.91"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.36%  '

$ws.Range('D29').Value = "'5.275"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.88%  '

$ws.Range('D30').Value = "'135.44"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.38%  '

$ws.Range('D31').Value = "'2.401"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.78%  '

$ws.Range('D32').Value = "'6.807"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.68%  '

$ws.Range('D33').Value = "'1.787.11"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.12%  '

$ws.Range('D34').Value = "'0.9579"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.75%  '

$ws.Range('D35').Value = "'0.02771"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.24%  '

$ws.Range('D36').Value = "'10.36"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.24%  '

$ws.Range('D37').Value = "'0.07373"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.94%  '

$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').Value = "'6.145"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.58%  '

$ws.Range('B39').Value = 'Algorand'
$ws.Range('C39').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D39').Value = "'0.2515"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.90%  '

$ws.Range('D40').Value = "'0.08779"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.29%  '

$ws.Range('D41').Value = "'1.396"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.53%  '

$ws.Range('D42').Value = "'0.7110"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.91%  '

$ws.Range('D43').Value = "'12.51"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.84%  '

$ws.Range('D44').Value = "'15.93"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +6.35%  '

$ws.Range('D45').Value = "'0.6540"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.57%  '

$ws.Range('D46').Value = "'2.336"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.02%  '

$ws.Range('B47').Value = 'PancakeSwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D47').Value = "'4.016"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.20%  '

$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').Value = "'133.04"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.97%  '

$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = "'0.07979"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.94%  '

$ws.Range('B50').Value = 'Flow'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range('D50').Value = "'1.201"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.55%  '

$ws.Range('B51').Value = 'ThetaToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D51').Value = "'1.195"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.16%  '
